# Auto-generated edit script: update price/profit figures across leve-profit sheets
# (mirrors a scheduled market-data refresh run)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2892.0942
$ws.Range("J17").Value = 2846.75
$ws.Range("L17").Value = 8540.25
$ws.Range("N17").Value = -8876.25
$ws.Range("H48").Value = 9000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H51").Value = 4773.36
$ws.Range("I51").Value = 5315.737
$ws.Range("J51").Value = 3055.8333
$ws.Range("K51").Value = 5315.737
$ws.Range("L51").Value = 3055.8333
$ws.Range("M51").Value = -4831.737
$ws.Range("N51").Value = -4023.8333
$ws.Range("H53").Value = 1008.9231
$ws.Range("I53").Value = 429.1111
$ws.Range("J53").Value = 2313.5
$ws.Range("K53").Value = 429.1111
$ws.Range("L53").Value = 2313.5
$ws.Range("M53").Value = 207.8889
$ws.Range("N53").Value = -3587.5
$ws.Range("H56").Value = 9000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H74").Value = 6602.375
$ws.Range("I74").Value = 6375.8667
$ws.Range("K74").Value = 6375.8667
$ws.Range("M74").Value = -5439.8667
$ws.Range("H77").Value = 6602.375
$ws.Range("I77").Value = 6375.8667
$ws.Range("K77").Value = 31879.3335
$ws.Range("M77").Value = -27199.3335
$ws.Range("H86").Value = 250003600
$ws.Range("I86").Value = 250002780
$ws.Range("K86").Value = 250002780
$ws.Range("M86").Value = -250001657
$ws.Range("H89").Value = 250003600
$ws.Range("I89").Value = 250002780
$ws.Range("K89").Value = 1250013900
$ws.Range("M89").Value = -1250008284
$ws.Range("H96").Value = 2581.1428
$ws.Range("I96").Value = 3208.4
$ws.Range("K96").Value = 9625.200000000001
$ws.Range("M96").Value = -8252.200000000001
$ws.Range("H137").Value = 2223256.5
$ws.Range("I137").Value = 763.3333
$ws.Range("J137").Value = 4274788.5
$ws.Range("K137").Value = 2289.9999
$ws.Range("L137").Value = 12824365.5
$ws.Range("M137").Value = 260.0001000000002
$ws.Range("N137").Value = -12829465.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 87068.39999999999
$ws.Range("J24").Value = 87068.39999999999
$ws.Range("L24").Value = 87068.39999999999
$ws.Range("N24").Value = -87816.39999999999
$ws.Range("H92").Value = 56455.285
$ws.Range("J92").Value = 55864.5
$ws.Range("L92").Value = 55864.5
$ws.Range("N92").Value = -60856.5
$ws.Range("H96").Value = 66275.39999999999
$ws.Range("J96").Value = 66275.39999999999
$ws.Range("L96").Value = 66275.39999999999
$ws.Range("N96").Value = -71767.39999999999
$ws.Range("H97").Value = 3469
$ws.Range("I97").Value = 1462.25
$ws.Range("K97").Value = 1462.25
$ws.Range("M97").Value = -966.25
$ws.Range("H100").Value = 87068.39999999999
$ws.Range("J100").Value = 87068.39999999999
$ws.Range("L100").Value = 87068.39999999999
$ws.Range("N100").Value = -89232.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 849.1667
$ws.Range("J22").Value = 1950
$ws.Range("L22").Value = 1950
$ws.Range("N22").Value = -2296
$ws.Range("H94").Value = 1848.6666
$ws.Range("I94").Value = 1622
$ws.Range("J94").Value = 1962
$ws.Range("K94").Value = 1622
$ws.Range("L94").Value = 1962
$ws.Range("M94").Value = -1171
$ws.Range("N94").Value = -2864
$ws.Range("H105").Value = 2698.1875
$ws.Range("I105").Value = 1803.8572
$ws.Range("K105").Value = 1803.8572
$ws.Range("M105").Value = -56.85719999999992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2820.625
$ws.Range("I58").Value = 2521.0244
$ws.Range("J58").Value = 4575.4287
$ws.Range("K58").Value = 2521.0244
$ws.Range("L58").Value = 4575.4287
$ws.Range("M58").Value = -2318.0244
$ws.Range("N58").Value = -4981.4287
$ws.Range("H122").Value = 3473.3572
$ws.Range("I122").Value = 2921.926
$ws.Range("J122").Value = 4465.933
$ws.Range("K122").Value = 8765.778
$ws.Range("L122").Value = 13397.799
$ws.Range("M122").Value = -6315.778
$ws.Range("N122").Value = -18297.799
$ws.Range("H132").Value = 4069.5293
$ws.Range("I132").Value = 4068.8
$ws.Range("K132").Value = 12206.4
$ws.Range("M132").Value = -9676.400000000001
$ws.Range("H136").Value = 2820.625
$ws.Range("I136").Value = 2521.0244
$ws.Range("J136").Value = 4575.4287
$ws.Range("K136").Value = 7563.073199999999
$ws.Range("L136").Value = 13726.2861
$ws.Range("M136").Value = -5013.073199999999
$ws.Range("N136").Value = -18826.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 833408.4399999999
$ws.Range("J12").Value = 1000084.6
$ws.Range("L12").Value = 3000253.8
$ws.Range("N12").Value = -3000599.8
$ws.Range("H48").Value = 1200
$ws.Range("I48").Value = 1200
$ws.Range("K48").Value = 3600
$ws.Range("M48").Value = -3350
$ws.Range("H110").Value = 4513.5
$ws.Range("H113").Value = 1436.8572
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1436.8572
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4310.571599999999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8650.571599999999
$ws.Range("H132").Value = 716660.4399999999
$ws.Range("I132").Value = 2057.125
$ws.Range("J132").Value = 1669464.9
$ws.Range("K132").Value = 18514.125
$ws.Range("L132").Value = 15025184.1
$ws.Range("M132").Value = -15984.125
$ws.Range("N132").Value = -15030244.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1329.5333
$ws.Range("I97").Value = 686.1818
$ws.Range("J97").Value = 3098.75
$ws.Range("K97").Value = 686.1818
$ws.Range("L97").Value = 3098.75
$ws.Range("M97").Value = -190.1818
$ws.Range("N97").Value = -4090.75
$ws.Range("H126").Value = 3722.923
$ws.Range("I126").Value = 3485
$ws.Range("K126").Value = 10455
$ws.Range("M126").Value = -7985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2315.8667
$ws.Range("I16").Value = 1715.2727
$ws.Range("K16").Value = 1715.2727
$ws.Range("M16").Value = -1545.2727
$ws.Range("H55").Value = 306.3846
$ws.Range("I55").Value = 269.375
$ws.Range("J55").Value = 365.6
$ws.Range("K55").Value = 269.375
$ws.Range("L55").Value = 365.6
$ws.Range("M55").Value = -96.375
$ws.Range("N55").Value = -711.6
$ws.Range("I61").Value = 1328.7222
$ws.Range("J61").Value = 2840
$ws.Range("K61").Value = 1328.7222
$ws.Range("L61").Value = 2840
$ws.Range("M61").Value = -1126.7222
$ws.Range("N61").Value = -3244
$ws.Range("I113").Value = 1328.7222
$ws.Range("J113").Value = 2840
$ws.Range("K113").Value = 1328.7222
$ws.Range("L113").Value = 2840
$ws.Range("M113").Value = 841.2778000000001
$ws.Range("N113").Value = -7180
$ws.Range("H132").Value = 4741.5483
$ws.Range("I132").Value = 4442.269
$ws.Range("K132").Value = 13326.807
$ws.Range("M132").Value = -10796.807
$ws.Range("H136").Value = 21249.625
$ws.Range("I136").Value = 24998
$ws.Range("J136").Value = 20714.143
$ws.Range("K136").Value = 74994
$ws.Range("L136").Value = 62142.429
$ws.Range("M136").Value = -72444
$ws.Range("N136").Value = -67242.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2827.2
$ws.Range("I81").Value = 2090.4614
$ws.Range("J81").Value = 4195.4287
$ws.Range("K81").Value = 4180.9228
$ws.Range("L81").Value = 8390.857400000001
$ws.Range("M81").Value = -3119.9228
$ws.Range("N81").Value = -10512.8574
$ws.Range("H84").Value = 2827.2
$ws.Range("I84").Value = 2090.4614
$ws.Range("J84").Value = 4195.4287
$ws.Range("K84").Value = 20904.614
$ws.Range("L84").Value = 41954.287
$ws.Range("M84").Value = -15600.614
$ws.Range("N84").Value = -52562.287
$ws.Range("H107").Value = 306.4737
$ws.Range("J107").Value = 342.9091
$ws.Range("L107").Value = 1028.7273
$ws.Range("N107").Value = -4868.7273
$ws.Range("H113").Value = 585.3333
$ws.Range("I113").Value = 444.58334
$ws.Range("K113").Value = 1333.75002
$ws.Range("M113").Value = 836.2499800000001
$ws.Range("H132").Value = 3525
$ws.Range("I132").Value = 3525
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10575
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8045
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 8377111.5
$ws.Range("I136").Value = 1714.0625
$ws.Range("K136").Value = 5142.1875
$ws.Range("M136").Value = -2592.1875
